# Update "想去人数" (F column) counts for several rows across sheets
# "展览" (Exhibitions), "演出" (Performances), "全部类型" (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2810
$ws1.Range("F7").Value = 3028
$ws1.Range("F18").Value = 9643
$ws1.Range("F22").Value = 7616
$ws1.Range("F23").Value = 12159
$ws1.Range("F30").Value = 244
$ws1.Range("F38").Value = 1167
$ws1.Range("F42").Value = 585

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F21").Value = 19
$ws2.Range("F22").Value = 1
$ws2.Range("F24").Value = 16
$ws2.Range("F27").Value = 8

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2810
$ws4.Range("F10").Value = 3028
$ws4.Range("F23").Value = 9643
$ws4.Range("F26").Value = 7616
$ws4.Range("F27").Value = 12159
$ws4.Range("F36").Value = 244
$ws4.Range("F44").Value = 19
$ws4.Range("F46").Value = 585
$ws4.Range("F47").Value = 16
